# Applies the edit described by the diff: appends two more
# "pages" of the same 7-person roster (rows 16-22 and 23-29)
# to Sheet1, each with a fresh per-row token in column J and an
# incremented column-K counter, and appends the matching 4 rows
# to Sheet2. Formulas are re-entered per row so every cell keeps
# a live formula (not a frozen value) exactly like the original
# rows 2-15.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Sheet1")
$ws2 = $wb.Worksheets.Item("Sheet2")

# Sheet1 row 16 (clone of row 2)
$ws1.Range("A16").Formula = '=PROPER(IFERROR(LEFT(C16,FIND(CHAR(46),C16)-1),C16))'
$ws1.Range("B16").Formula = '=IFERROR(PROPER(RIGHT(C16,LEN(C16)-FIND("@",SUBSTITUTE(C16,".","@",((LEN(C16)-LEN(SUBSTITUTE(C16,".","")))/LEN("\")))))), "Unknown")'
$ws1.Range("C16").Formula = '=SUBSTITUTE(SUBSTITUTE(LOWER(LEFT(D16,FIND(CHAR(64),D16)-1)),CHAR(45),CHAR(46)),CHAR(95),CHAR(46))'
$ws1.Range("D16").Value = 'sourabh.awasthi@capgemini.com'
$ws1.Range("E16").Formula = '=LEFT(H16,FIND(CHAR(46),H16)-1)'
$ws1.Range("F16").Formula = '=CONCATENATE("ITPartner\",I16)'
$ws1.Range("H16").Formula = '=RIGHT(D16,LEN(D16)-FIND(CHAR(64),D16))'
$ws1.Range("I16").Formula = '=PROPER(E16)'
$ws1.Range("J16").Value = '%+ZN#w1)4TAa'
$ws1.Range("K16").Value = 83
$ws1.Range("M16").Value = $true
$ws1.Range("P16").Formula = '=COUNTIF(D:D,D16)'

# Sheet1 row 17 (clone of row 3)
$ws1.Range("A17").Formula = '=PROPER(IFERROR(LEFT(C17,FIND(CHAR(46),C17)-1),C17))'
$ws1.Range("B17").Formula = '=IFERROR(PROPER(RIGHT(C17,LEN(C17)-FIND("@",SUBSTITUTE(C17,".","@",((LEN(C17)-LEN(SUBSTITUTE(C17,".","")))/LEN("\")))))), "Unknown")'
$ws1.Range("C17").Formula = '=SUBSTITUTE(SUBSTITUTE(LOWER(LEFT(D17,FIND(CHAR(64),D17)-1)),CHAR(45),CHAR(46)),CHAR(95),CHAR(46))'
$ws1.Range("D17").Value = 'sandipan.deb@capgemini.com'
$ws1.Range("E17").Formula = '=LEFT(H17,FIND(CHAR(46),H17)-1)'
$ws1.Range("F17").Formula = '=CONCATENATE("ITPartner\",I17)'
$ws1.Range("H17").Formula = '=RIGHT(D17,LEN(D17)-FIND(CHAR(64),D17))'
$ws1.Range("I17").Formula = '=PROPER(E17)'
$ws1.Range("J17").Value = '"6UUl0IJJ(L['
$ws1.Range("K17").Value = 83
$ws1.Range("M17").Value = $true
$ws1.Range("P17").Formula = '=COUNTIF(D:D,D17)'

# Sheet1 row 18 (clone of row 4)
$ws1.Range("A18").Formula = '=PROPER(IFERROR(LEFT(C18,FIND(CHAR(46),C18)-1),C18))'
$ws1.Range("B18").Formula = '=IFERROR(PROPER(RIGHT(C18,LEN(C18)-FIND("@",SUBSTITUTE(C18,".","@",((LEN(C18)-LEN(SUBSTITUTE(C18,".","")))/LEN("\")))))), "Unknown")'
$ws1.Range("C18").Formula = '=SUBSTITUTE(SUBSTITUTE(LOWER(LEFT(D18,FIND(CHAR(64),D18)-1)),CHAR(45),CHAR(46)),CHAR(95),CHAR(46))'
$ws1.Range("D18").Value = 'biswaji.deb@capgemini.com'
$ws1.Range("E18").Formula = '=LEFT(H18,FIND(CHAR(46),H18)-1)'
$ws1.Range("F18").Formula = '=CONCATENATE("ITPartner\",I18)'
$ws1.Range("H18").Formula = '=RIGHT(D18,LEN(D18)-FIND(CHAR(64),D18))'
$ws1.Range("I18").Formula = '=PROPER(E18)'
$ws1.Range("J18").Value = 'L09Pw&yUfcyl'
$ws1.Range("K18").Value = 83
$ws1.Range("M18").Value = $true
$ws1.Range("P18").Formula = '=COUNTIF(D:D,D18)'

# Sheet1 row 19 (clone of row 5)
$ws1.Range("A19").Formula = '=PROPER(IFERROR(LEFT(C19,FIND(CHAR(46),C19)-1),C19))'
$ws1.Range("B19").Formula = '=IFERROR(PROPER(RIGHT(C19,LEN(C19)-FIND("@",SUBSTITUTE(C19,".","@",((LEN(C19)-LEN(SUBSTITUTE(C19,".","")))/LEN("\")))))), "Unknown")'
$ws1.Range("C19").Formula = '=SUBSTITUTE(SUBSTITUTE(LOWER(LEFT(D19,FIND(CHAR(64),D19)-1)),CHAR(45),CHAR(46)),CHAR(95),CHAR(46))'
$ws1.Range("D19").Value = 'debanjan.das@capgemini.com'
$ws1.Range("E19").Formula = '=LEFT(H19,FIND(CHAR(46),H19)-1)'
$ws1.Range("F19").Formula = '=CONCATENATE("ITPartner\",I19)'
$ws1.Range("H19").Formula = '=RIGHT(D19,LEN(D19)-FIND(CHAR(64),D19))'
$ws1.Range("I19").Formula = '=PROPER(E19)'
$ws1.Range("J19").Value = 'a$3c&C6e/xF5'
$ws1.Range("K19").Value = 83
$ws1.Range("M19").Value = $true
$ws1.Range("P19").Formula = '=COUNTIF(D:D,D19)'

# Sheet1 row 20 (clone of row 6)
$ws1.Range("A20").Formula = '=PROPER(IFERROR(LEFT(C20,FIND(CHAR(46),C20)-1),C20))'
$ws1.Range("B20").Formula = '=IFERROR(PROPER(RIGHT(C20,LEN(C20)-FIND("@",SUBSTITUTE(C20,".","@",((LEN(C20)-LEN(SUBSTITUTE(C20,".","")))/LEN("\")))))), "Unknown")'
$ws1.Range("C20").Formula = '=SUBSTITUTE(SUBSTITUTE(LOWER(LEFT(D20,FIND(CHAR(64),D20)-1)),CHAR(45),CHAR(46)),CHAR(95),CHAR(46))'
$ws1.Range("D20").Value = 'dhiraj.kajari@capgemini.com'
$ws1.Range("E20").Formula = '=LEFT(H20,FIND(CHAR(46),H20)-1)'
$ws1.Range("F20").Formula = '=CONCATENATE("ITPartner\",I20)'
$ws1.Range("H20").Formula = '=RIGHT(D20,LEN(D20)-FIND(CHAR(64),D20))'
$ws1.Range("I20").Formula = '=PROPER(E20)'
$ws1.Range("J20").Value = 'T5IVsi6cu2*6'
$ws1.Range("K20").Value = 83
$ws1.Range("M20").Value = $true
$ws1.Range("P20").Formula = '=COUNTIF(D:D,D20)'

# Sheet1 row 21 (clone of row 7)
$ws1.Range("A21").Formula = '=PROPER(IFERROR(LEFT(C21,FIND(CHAR(46),C21)-1),C21))'
$ws1.Range("B21").Formula = '=IFERROR(PROPER(RIGHT(C21,LEN(C21)-FIND("@",SUBSTITUTE(C21,".","@",((LEN(C21)-LEN(SUBSTITUTE(C21,".","")))/LEN("\")))))), "Unknown")'
$ws1.Range("C21").Formula = '=SUBSTITUTE(SUBSTITUTE(LOWER(LEFT(D21,FIND(CHAR(64),D21)-1)),CHAR(45),CHAR(46)),CHAR(95),CHAR(46))'
$ws1.Range("D21").Value = 'manoj-kumar.b.s@capgemini.com'
$ws1.Range("E21").Formula = '=LEFT(H21,FIND(CHAR(46),H21)-1)'
$ws1.Range("F21").Formula = '=CONCATENATE("ITPartner\",I21)'
$ws1.Range("H21").Formula = '=RIGHT(D21,LEN(D21)-FIND(CHAR(64),D21))'
$ws1.Range("I21").Formula = '=PROPER(E21)'
$ws1.Range("J21").Value = 'v%/&NL(l}96b'
$ws1.Range("K21").Value = 83
$ws1.Range("M21").Value = $true
$ws1.Range("P21").Formula = '=COUNTIF(D:D,D21)'

# Sheet1 row 22 (clone of row 8)
$ws1.Range("A22").Formula = '=PROPER(IFERROR(LEFT(C22,FIND(CHAR(46),C22)-1),C22))'
$ws1.Range("B22").Formula = '=IFERROR(PROPER(RIGHT(C22,LEN(C22)-FIND("@",SUBSTITUTE(C22,".","@",((LEN(C22)-LEN(SUBSTITUTE(C22,".","")))/LEN("\")))))), "Unknown")'
$ws1.Range("C22").Formula = '=SUBSTITUTE(SUBSTITUTE(LOWER(LEFT(D22,FIND(CHAR(64),D22)-1)),CHAR(45),CHAR(46)),CHAR(95),CHAR(46))'
$ws1.Range("D22").Value = 'mayur.bhorkar@capgemini.com'
$ws1.Range("E22").Formula = '=LEFT(H22,FIND(CHAR(46),H22)-1)'
$ws1.Range("F22").Formula = '=CONCATENATE("ITPartner\",I22)'
$ws1.Range("H22").Formula = '=RIGHT(D22,LEN(D22)-FIND(CHAR(64),D22))'
$ws1.Range("I22").Formula = '=PROPER(E22)'
$ws1.Range("J22").Value = '!SldB/c8Tc&x'
$ws1.Range("K22").Value = 83
$ws1.Range("M22").Value = $true
$ws1.Range("P22").Formula = '=COUNTIF(D:D,D22)'

# Sheet1 row 23 (clone of row 2)
$ws1.Range("A23").Formula = '=PROPER(IFERROR(LEFT(C23,FIND(CHAR(46),C23)-1),C23))'
$ws1.Range("B23").Formula = '=IFERROR(PROPER(RIGHT(C23,LEN(C23)-FIND("@",SUBSTITUTE(C23,".","@",((LEN(C23)-LEN(SUBSTITUTE(C23,".","")))/LEN("\")))))), "Unknown")'
$ws1.Range("C23").Formula = '=SUBSTITUTE(SUBSTITUTE(LOWER(LEFT(D23,FIND(CHAR(64),D23)-1)),CHAR(45),CHAR(46)),CHAR(95),CHAR(46))'
$ws1.Range("D23").Value = 'sourabh.awasthi@capgemini.com'
$ws1.Range("E23").Formula = '=LEFT(H23,FIND(CHAR(46),H23)-1)'
$ws1.Range("F23").Formula = '=CONCATENATE("ITPartner\",I23)'
$ws1.Range("H23").Formula = '=RIGHT(D23,LEN(D23)-FIND(CHAR(64),D23))'
$ws1.Range("I23").Formula = '=PROPER(E23)'
$ws1.Range("J23").Value = '"iRC#%@GY[Dw'
$ws1.Range("K23").Value = 84
$ws1.Range("M23").Value = $true
$ws1.Range("P23").Formula = '=COUNTIF(D:D,D23)'

# Sheet1 row 24 (clone of row 3)
$ws1.Range("A24").Formula = '=PROPER(IFERROR(LEFT(C24,FIND(CHAR(46),C24)-1),C24))'
$ws1.Range("B24").Formula = '=IFERROR(PROPER(RIGHT(C24,LEN(C24)-FIND("@",SUBSTITUTE(C24,".","@",((LEN(C24)-LEN(SUBSTITUTE(C24,".","")))/LEN("\")))))), "Unknown")'
$ws1.Range("C24").Formula = '=SUBSTITUTE(SUBSTITUTE(LOWER(LEFT(D24,FIND(CHAR(64),D24)-1)),CHAR(45),CHAR(46)),CHAR(95),CHAR(46))'
$ws1.Range("D24").Value = 'sandipan.deb@capgemini.com'
$ws1.Range("E24").Formula = '=LEFT(H24,FIND(CHAR(46),H24)-1)'
$ws1.Range("F24").Formula = '=CONCATENATE("ITPartner\",I24)'
$ws1.Range("H24").Formula = '=RIGHT(D24,LEN(D24)-FIND(CHAR(64),D24))'
$ws1.Range("I24").Formula = '=PROPER(E24)'
$ws1.Range("J24").Value = 'QW+*EBu9Aysv'
$ws1.Range("K24").Value = 84
$ws1.Range("M24").Value = $true
$ws1.Range("P24").Formula = '=COUNTIF(D:D,D24)'

# Sheet1 row 25 (clone of row 4)
$ws1.Range("A25").Formula = '=PROPER(IFERROR(LEFT(C25,FIND(CHAR(46),C25)-1),C25))'
$ws1.Range("B25").Formula = '=IFERROR(PROPER(RIGHT(C25,LEN(C25)-FIND("@",SUBSTITUTE(C25,".","@",((LEN(C25)-LEN(SUBSTITUTE(C25,".","")))/LEN("\")))))), "Unknown")'
$ws1.Range("C25").Formula = '=SUBSTITUTE(SUBSTITUTE(LOWER(LEFT(D25,FIND(CHAR(64),D25)-1)),CHAR(45),CHAR(46)),CHAR(95),CHAR(46))'
$ws1.Range("D25").Value = 'biswaji.deb@capgemini.com'
$ws1.Range("E25").Formula = '=LEFT(H25,FIND(CHAR(46),H25)-1)'
$ws1.Range("F25").Formula = '=CONCATENATE("ITPartner\",I25)'
$ws1.Range("H25").Formula = '=RIGHT(D25,LEN(D25)-FIND(CHAR(64),D25))'
$ws1.Range("I25").Formula = '=PROPER(E25)'
$ws1.Range("J25").Value = '}]VDNkYX/k7{'
$ws1.Range("K25").Value = 84
$ws1.Range("M25").Value = $true
$ws1.Range("P25").Formula = '=COUNTIF(D:D,D25)'

# Sheet1 row 26 (clone of row 5)
$ws1.Range("A26").Formula = '=PROPER(IFERROR(LEFT(C26,FIND(CHAR(46),C26)-1),C26))'
$ws1.Range("B26").Formula = '=IFERROR(PROPER(RIGHT(C26,LEN(C26)-FIND("@",SUBSTITUTE(C26,".","@",((LEN(C26)-LEN(SUBSTITUTE(C26,".","")))/LEN("\")))))), "Unknown")'
$ws1.Range("C26").Formula = '=SUBSTITUTE(SUBSTITUTE(LOWER(LEFT(D26,FIND(CHAR(64),D26)-1)),CHAR(45),CHAR(46)),CHAR(95),CHAR(46))'
$ws1.Range("D26").Value = 'debanjan.das@capgemini.com'
$ws1.Range("E26").Formula = '=LEFT(H26,FIND(CHAR(46),H26)-1)'
$ws1.Range("F26").Formula = '=CONCATENATE("ITPartner\",I26)'
$ws1.Range("H26").Formula = '=RIGHT(D26,LEN(D26)-FIND(CHAR(64),D26))'
$ws1.Range("I26").Formula = '=PROPER(E26)'
$ws1.Range("J26").Value = 'z7Q[k7+3{}5a'
$ws1.Range("K26").Value = 84
$ws1.Range("M26").Value = $true
$ws1.Range("P26").Formula = '=COUNTIF(D:D,D26)'

# Sheet1 row 27 (clone of row 6)
$ws1.Range("A27").Formula = '=PROPER(IFERROR(LEFT(C27,FIND(CHAR(46),C27)-1),C27))'
$ws1.Range("B27").Formula = '=IFERROR(PROPER(RIGHT(C27,LEN(C27)-FIND("@",SUBSTITUTE(C27,".","@",((LEN(C27)-LEN(SUBSTITUTE(C27,".","")))/LEN("\")))))), "Unknown")'
$ws1.Range("C27").Formula = '=SUBSTITUTE(SUBSTITUTE(LOWER(LEFT(D27,FIND(CHAR(64),D27)-1)),CHAR(45),CHAR(46)),CHAR(95),CHAR(46))'
$ws1.Range("D27").Value = 'dhiraj.kajari@capgemini.com'
$ws1.Range("E27").Formula = '=LEFT(H27,FIND(CHAR(46),H27)-1)'
$ws1.Range("F27").Formula = '=CONCATENATE("ITPartner\",I27)'
$ws1.Range("H27").Formula = '=RIGHT(D27,LEN(D27)-FIND(CHAR(64),D27))'
$ws1.Range("I27").Formula = '=PROPER(E27)'
$ws1.Range("J27").Value = 'p=fCvU}BYlAA'
$ws1.Range("K27").Value = 84
$ws1.Range("M27").Value = $true
$ws1.Range("P27").Formula = '=COUNTIF(D:D,D27)'

# Sheet1 row 28 (clone of row 7)
$ws1.Range("A28").Formula = '=PROPER(IFERROR(LEFT(C28,FIND(CHAR(46),C28)-1),C28))'
$ws1.Range("B28").Formula = '=IFERROR(PROPER(RIGHT(C28,LEN(C28)-FIND("@",SUBSTITUTE(C28,".","@",((LEN(C28)-LEN(SUBSTITUTE(C28,".","")))/LEN("\")))))), "Unknown")'
$ws1.Range("C28").Formula = '=SUBSTITUTE(SUBSTITUTE(LOWER(LEFT(D28,FIND(CHAR(64),D28)-1)),CHAR(45),CHAR(46)),CHAR(95),CHAR(46))'
$ws1.Range("D28").Value = 'manoj-kumar.b.s@capgemini.com'
$ws1.Range("E28").Formula = '=LEFT(H28,FIND(CHAR(46),H28)-1)'
$ws1.Range("F28").Formula = '=CONCATENATE("ITPartner\",I28)'
$ws1.Range("H28").Formula = '=RIGHT(D28,LEN(D28)-FIND(CHAR(64),D28))'
$ws1.Range("I28").Formula = '=PROPER(E28)'
$ws1.Range("J28").Value = 'yw?5Cu{vPZel'
$ws1.Range("K28").Value = 84
$ws1.Range("M28").Value = $true
$ws1.Range("P28").Formula = '=COUNTIF(D:D,D28)'

# Sheet1 row 29 (clone of row 8)
$ws1.Range("A29").Formula = '=PROPER(IFERROR(LEFT(C29,FIND(CHAR(46),C29)-1),C29))'
$ws1.Range("B29").Formula = '=IFERROR(PROPER(RIGHT(C29,LEN(C29)-FIND("@",SUBSTITUTE(C29,".","@",((LEN(C29)-LEN(SUBSTITUTE(C29,".","")))/LEN("\")))))), "Unknown")'
$ws1.Range("C29").Formula = '=SUBSTITUTE(SUBSTITUTE(LOWER(LEFT(D29,FIND(CHAR(64),D29)-1)),CHAR(45),CHAR(46)),CHAR(95),CHAR(46))'
$ws1.Range("D29").Value = 'mayur.bhorkar@capgemini.com'
$ws1.Range("E29").Formula = '=LEFT(H29,FIND(CHAR(46),H29)-1)'
$ws1.Range("F29").Formula = '=CONCATENATE("ITPartner\",I29)'
$ws1.Range("H29").Formula = '=RIGHT(D29,LEN(D29)-FIND(CHAR(64),D29))'
$ws1.Range("I29").Formula = '=PROPER(E29)'
$ws1.Range("J29").Value = '@pvZlr*!RWL&'
$ws1.Range("K29").Value = 84
$ws1.Range("M29").Value = $true
$ws1.Range("P29").Formula = '=COUNTIF(D:D,D29)'

# Sheet2 rows 6-9 (clone of rows 2-5)
$ws2.Range("C6").Value = 'sandipan.deb'
$ws2.Range("I6").Value = 'Capgemini'
$ws2.Range("C7").Value = 'mayur.bhorkar'
$ws2.Range("I7").Value = 'Capgemini'
$ws2.Range("C8").Value = 'sandipan.deb'
$ws2.Range("I8").Value = 'Capgemini'
$ws2.Range("C9").Value = 'mayur.bhorkar'
$ws2.Range("I9").Value = 'Capgemini'

$wb.Application.Calculate()
